$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Fullname" validation column (M) to let callers sanity-check
# that the fullname field doesn't contain digits.
$ws.Range("M1").Value = "Fullname"

$ws.Range("M2").Value = "Bach Hoang"
$ws.Range("M3").Value = "Bach Hoang"
$ws.Range("M4").Value = "Bach Hoang"
$ws.Range("M5").Value = "Bach Hoang"

$ws.Range("M4").Select()
